$d = $word.ActiveDocument

function ReplaceInRange($range, $findText, $replaceText) {
    $r = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $r) {
        Write-Host "WARNING: could not find '$findText'"
    }
    return $r
}

# --- Paragraph 2: title "固件升级" -> "Firmware Upgrade" ---
$p = $d.Paragraphs.Item(2)
ReplaceInRange $p.Range "固件升级" "Firmware Upgrade"

# --- Paragraph 3: intro sentence ---
$p = $d.Paragraphs.Item(3)
ReplaceInRange $p.Range "使用" "Use"
ReplaceInRange $p.Range "Fly App或者DJI Assistant 2(Cons" "Fly App or DJI Assistant 2 (Cons"
ReplaceInRange $p.Range "umer Drones Series)调参软件对飞行器和遥控器" "umer Drones Series) to update aircraft and remote"
ReplaceInRange $p.Range "进行升级。" "."

# --- Paragraph 5: bold heading "使用DJI Fly App升 级" ---
$p = $d.Paragraphs.Item(5)
ReplaceInRange $p.Range "使用" "Use"
ReplaceInRange $p.Range "升" "litter"
ReplaceInRange $p.Range "级" "level"
